# Generate Report for Handback
# d001c143-1927-4d3e-a149-6e447a561c0d and eb642c3d-0a35-4ceb-9ab7-ebde167e0110
# have now been handed back (in sync with en-US) for both locales, so the
# Overview sheet's status columns, and each locale sheet's Status / Latest
# Target File / Latest Handback File / Latest Handback DateTime columns need
# to be updated to reflect that.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (d001c143...) and 5 (eb642c3d...) move from
# "Ready for handoff" to "Handed back: in sync with en-US" in both the
# zh-cn and de-de status columns (B and C).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B4").Value = $handedBack
$wsOverview.Range("C4").Value = $handedBack
$wsOverview.Range("B5").Value = $handedBack
$wsOverview.Range("C5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 - d001c143-1927-4d3e-a149-6e447a561c0d
$wsZh.Range("C4").Value = $handedBack
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/46eae1539aceb462f9dc3f8d558506ddc4f3b71f/e2e/d001c143-1927-4d3e-a149-6e447a561c0d.md", [System.Type]::Missing, [System.Type]::Missing, "d001c143-1927-4d3e-a149-6e447a561c0d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5905dbad2ed5cb1459e44b3bdb7a9e232e2498ba/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/d001c143-1927-4d3e-a149-6e447a561c0d.5b9cb1e5e368ed14608e4ed97f2982de493b081c.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "d001c143-1927-4d3e-a149-6e447a561c0d.5b9cb1e5e368ed14608e4ed97f2982de493b081c.zh-cn.xlf")
$wsZh.Range("H4").Value = "2016-03-17 03:48:00"

# Row 5 - eb642c3d-0a35-4ceb-9ab7-ebde167e0110
$wsZh.Range("C5").Value = $handedBack
$wsZh.Hyperlinks.Add($wsZh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/46eae1539aceb462f9dc3f8d558506ddc4f3b71f/e2e/eb642c3d-0a35-4ceb-9ab7-ebde167e0110.md", [System.Type]::Missing, [System.Type]::Missing, "eb642c3d-0a35-4ceb-9ab7-ebde167e0110.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5905dbad2ed5cb1459e44b3bdb7a9e232e2498ba/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/eb642c3d-0a35-4ceb-9ab7-ebde167e0110.9b821621b6c36d021dbcaf79bd9a8c1c712758bb.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "eb642c3d-0a35-4ceb-9ab7-ebde167e0110.9b821621b6c36d021dbcaf79bd9a8c1c712758bb.zh-cn.xlf")
$wsZh.Range("H5").Value = "2016-03-17 03:48:00"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 - d001c143-1927-4d3e-a149-6e447a561c0d
$wsDe.Range("C4").Value = $handedBack
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b3df44fd3212a4f6a88544a767f8f0a9b8359dde/e2e/d001c143-1927-4d3e-a149-6e447a561c0d.md", [System.Type]::Missing, [System.Type]::Missing, "d001c143-1927-4d3e-a149-6e447a561c0d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a8aa392b019d6725dee04bb172f42f5407658a4a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/d001c143-1927-4d3e-a149-6e447a561c0d.5b9cb1e5e368ed14608e4ed97f2982de493b081c.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "d001c143-1927-4d3e-a149-6e447a561c0d.5b9cb1e5e368ed14608e4ed97f2982de493b081c.de-de.xlf")
$wsDe.Range("H4").Value = "2016-03-17 03:48:13"

# Row 5 - eb642c3d-0a35-4ceb-9ab7-ebde167e0110
$wsDe.Range("C5").Value = $handedBack
$wsDe.Hyperlinks.Add($wsDe.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b3df44fd3212a4f6a88544a767f8f0a9b8359dde/e2e/eb642c3d-0a35-4ceb-9ab7-ebde167e0110.md", [System.Type]::Missing, [System.Type]::Missing, "eb642c3d-0a35-4ceb-9ab7-ebde167e0110.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a8aa392b019d6725dee04bb172f42f5407658a4a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/eb642c3d-0a35-4ceb-9ab7-ebde167e0110.9b821621b6c36d021dbcaf79bd9a8c1c712758bb.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "eb642c3d-0a35-4ceb-9ab7-ebde167e0110.9b821621b6c36d021dbcaf79bd9a8c1c712758bb.de-de.xlf")
$wsDe.Range("H5").Value = "2016-03-17 03:48:13"
